$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Floor-type lookup table (G1:H4, already an Excel Table "Table2") ---
$ws.Range("G2").Value = "Hardwood"
$ws.Range("H2").Value = 1.39

$ws.Range("G3").Value = "Carpet"
$ws.Range("H3").Value = 3.99

$ws.Range("G4").Value = "Tile"
$ws.Range("H4").Value = 4.99

# --- Test-case rows (A3:D7) ---
$ws.Range("A3").Value = "Hardwood"
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 20
$ws.Range("D3").Formula = "=(B3*C3)*H2"

$ws.Range("A4").Value = "Hardwood"
$ws.Range("B4").Value = 20.5
$ws.Range("C4").Value = 15.2
$ws.Range("D4").Formula = "=(B4*C4)*H2"

$ws.Range("A5").Value = "Carpet"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 15
$ws.Range("D5").Formula = "=(B5*C5)*H3"

$ws.Range("A6").Value = "Carpet"
$ws.Range("B6").Value = 14.239
$ws.Range("C6").Value = 10.4
$ws.Range("D6").Formula = "=(B6*C6)*H3"

$ws.Range("A7").Value = "Tile"
$ws.Range("B7").Value = 16.05
$ws.Range("C7").Value = 10
$ws.Range("D7").Formula = "=(B7*C7)*H4"

# Currency number format for the computed room-cost column
$ws.Range("D3:D7").NumberFormat = """$""#,##0.00"

# Final cursor position left by the author
$ws.Range("C7").Select() | Out-Null
